# Update gh-pages output data (想去人数 / 最低票价 columns) for the
# 展览, 演出 and 全部类型 sheets, matching the freshly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1694
$ws.Range("G2").Value = 68
$ws.Range("F3").Value = 9191
$ws.Range("G3").Value = 128
$ws.Range("F4").Value = 119
$ws.Range("F5").Value = 510
$ws.Range("F7").Value = 1386
$ws.Range("F9").Value = 66
$ws.Range("F10").Value = 102
$ws.Range("F11").Value = 5970
$ws.Range("F13").Value = 393
$ws.Range("F14").Value = 106
$ws.Range("F15").Value = 4651
$ws.Range("F19").Value = 38
$ws.Range("F20").Value = 345
$ws.Range("F25").Value = 3071

# --- 演出 (Performance) sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 49

# --- 全部类型 (All types) sheet ----------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1694
$ws.Range("G2").Value = 68
$ws.Range("F3").Value = 9191
$ws.Range("G3").Value = 128
$ws.Range("F4").Value = 119
$ws.Range("F5").Value = 49
$ws.Range("F6").Value = 510
$ws.Range("F8").Value = 1386
$ws.Range("F10").Value = 66
$ws.Range("F11").Value = 102
$ws.Range("F12").Value = 5970
$ws.Range("F14").Value = 393
$ws.Range("F15").Value = 106
$ws.Range("F16").Value = 4651
$ws.Range("F20").Value = 38
$ws.Range("F21").Value = 345
$ws.Range("F26").Value = 3071
